$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.134.68"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.835.33"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("D4").Formula = "'0.9991"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Formula = "'240.24"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").Formula = "'0.6649"
$ws.Range("E6").Value = "  -4.48%  "
$ws.Range("D7").Formula = "'1.000"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Formula = "'0.2955"
$ws.Range("E8").Value = "  -3.94%  "
$ws.Range("D9").Formula = "'0.07360"
$ws.Range("E9").Value = "  -4.39%  "
$ws.Range("D10").Formula = "'22.77"
$ws.Range("E10").Value = "  -3.66%  "
$ws.Range("D11").Formula = "'0.07682"
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "1.836.98"
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Formula = "'5.020"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").Formula = "'0.6751"
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Formula = "'86.23"
$ws.Range("E15").Value = "  -5.48%  "
$ws.Range("D16").Formula = "'6.183"
$ws.Range("E16").Value = "  -2.52%  "
$ws.Range("D17").Value = "29.055.02"
$ws.Range("E17").Value = "  -1.45%  "
$ws.Range("D18").Formula = "'0.000008241"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").Formula = "'228.82"
$ws.Range("E19").Value = "  -4.24%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("D21").Formula = "'0.9996"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Formula = "'7.299"
$ws.Range("E22").Value = "  -4.24%  "
$ws.Range("D23").Formula = "'1.0000"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Formula = "'161.05"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Formula = "'0.1418"
$ws.Range("E25").Value = "  -5.25%  "
$ws.Range("D26").Formula = "'8.681"
$ws.Range("E26").Value = "  -2.38%  "
$ws.Range("D27").Formula = "'18.03"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Formula = "'1.503"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Formula = "'4.232"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").Formula = "'4.101"
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("D31").Formula = "'1.203"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Formula = "'0.05310"
$ws.Range("E32").Value = "  +3.90%  "
$ws.Range("D33").Formula = "'1.859"
$ws.Range("D34").Formula = "'0.7467"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "1.316.69"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("E38").Value = "  -3.87%  "
$ws.Range("D40").Formula = "'0.9226"
$ws.Range("E40").Value = "  -3.27%  "
$ws.Range("D41").Formula = "'5.987"
$ws.Range("E41").Value = "  +3.68%  "
$ws.Range("D42").Formula = "'0.9989"
$ws.Range("E42").Value = "  -0.19%  "
$ws.Range("D43").Formula = "'103.41"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "1.985.46"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("D45").Formula = "'0.5170"
$ws.Range("E45").Value = "  -1.24%  "
$ws.Range("E46").Value = "  -3.13%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Formula = "'63.73"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Formula = "'1.761"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").Formula = "'9.289"
$ws.Range("E49").Value = "  -5.62%  "
$ws.Range("D50").Formula = "'0.07468"
$ws.Range("E50").Value = "  +9.39%  "
$ws.Range("D51").Formula = "'0.05930"
